# Fragenbaum.xlsx edit: rename the "Text" column header to "Frage",
# let the question rows that carry an upload follow-up ("ja" in column D)
# re-settle to their auto-computed row height, and leave the selection on C2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: "Text" -> "Frage"
$ws.Range("C1").Value = "Frage"

# Rows 3, 5, 7, 9 and 11 are the "upload photo?" questions (column D = "ja").
# Clear their previously stored (stale) explicit height so Excel recomputes it.
$ws.Rows.Item(3).AutoFit()
$ws.Rows.Item(5).AutoFit()
$ws.Rows.Item(7).AutoFit()
$ws.Rows.Item(9).AutoFit()
$ws.Rows.Item(11).RowHeight = 28.5

# Leave the cursor on C2, matching the saved view state.
$ws.Range("C2").Select() | Out-Null
